$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates: ref -> new value. Values that look like plain numbers in the
# "Price" column must stay TEXT (the sheet stores prices as literal strings,
# e.g. "30.617.15" / "1.000", so trailing zeros and dotted groupings survive).
$updates = @(
    @{ Cell = "D2"; Value = '30.617.15' }
    @{ Cell = "E2"; Value = '  +0.82%  ' }
    @{ Cell = "D3"; Value = '1.878.20' }
    @{ Cell = "D4"; Value = '1.000' }
    @{ Cell = "E4"; Value = '  -0.02%  ' }
    @{ Cell = "D5"; Value = '248.20' }
    @{ Cell = "E5"; Value = '  +1.50%  ' }
    @{ Cell = "D7"; Value = '0.4761' }
    @{ Cell = "E7"; Value = '  -0.08%  ' }
    @{ Cell = "D8"; Value = '0.2925' }
    @{ Cell = "E8"; Value = '  +1.71%  ' }
    @{ Cell = "E9"; Value = '  +0.15%  ' }
    @{ Cell = "D10"; Value = '22.00' }
    @{ Cell = "E10"; Value = '  +3.40%  ' }
    @{ Cell = "D11"; Value = '0.07737' }
    @{ Cell = "E11"; Value = '  -0.21%  ' }
    @{ Cell = "D12"; Value = '0.7403' }
    @{ Cell = "E12"; Value = '  +0.84%  ' }
    @{ Cell = "D13"; Value = '96.80' }
    @{ Cell = "E13"; Value = '  +0.21%  ' }
    @{ Cell = "D14"; Value = '1.877.12' }
    @{ Cell = "E14"; Value = '  -0.20%  ' }
    @{ Cell = "E15"; Value = '  +1.69%  ' }
    @{ Cell = "D16"; Value = '274.22' }
    @{ Cell = "E16"; Value = '  +0.22%  ' }
    @{ Cell = "D17"; Value = '30.720.53' }
    @{ Cell = "E17"; Value = '  +1.20%  ' }
    @{ Cell = "D18"; Value = '13.25' }
    @{ Cell = "E18"; Value = '  -0.93%  ' }
    @{ Cell = "D19"; Value = '0.000007531' }
    @{ Cell = "E19"; Value = '  -0.07%  ' }
    @{ Cell = "D20"; Value = '0.9998' }
    @{ Cell = "D21"; Value = '2.121.94' }
    @{ Cell = "E21"; Value = '  -0.27%  ' }
    @{ Cell = "D22"; Value = '1.001' }
    @{ Cell = "E22"; Value = '  -0.03%  ' }
    @{ Cell = "D23"; Value = '5.256' }
    @{ Cell = "E23"; Value = '  +0.52%  ' }
    @{ Cell = "E24"; Value = '  +0.71%  ' }
    @{ Cell = "D25"; Value = '9.200' }
    @{ Cell = "E25"; Value = '  -0.49%  ' }
    @{ Cell = "E26"; Value = '  +1.26%  ' }
    @{ Cell = "D27"; Value = '18.88' }
    @{ Cell = "E27"; Value = '  -0.18%  ' }
    @{ Cell = "D28"; Value = '1.916' }
    @{ Cell = "E28"; Value = '  -2.06%  ' }
    @{ Cell = "D29"; Value = '0.09850' }
    @{ Cell = "E29"; Value = '  -1.24%  ' }
    @{ Cell = "E30"; Value = '  -2.42%  ' }
    @{ Cell = "D31"; Value = '1.503' }
    @{ Cell = "E31"; Value = '  -0.17%  ' }
    @{ Cell = "D32"; Value = '4.290' }
    @{ Cell = "E32"; Value = '  -0.40%  ' }
    @{ Cell = "D33"; Value = '4.117' }
    @{ Cell = "E33"; Value = '  +1.10%  ' }
    @{ Cell = "D34"; Value = '0.04837' }
    @{ Cell = "E34"; Value = '  +2.00%  ' }
    @{ Cell = "D35"; Value = '1.126' }
    @{ Cell = "E35"; Value = '  +0.46%  ' }
    @{ Cell = "D36"; Value = '0.6960' }
    @{ Cell = "E36"; Value = '  +0.14%  ' }
    @{ Cell = "E37"; Value = '  +0.00%  ' }
    @{ Cell = "D38"; Value = '0.01879' }
    @{ Cell = "E38"; Value = '  +1.33%  ' }
    @{ Cell = "D39"; Value = '2.764' }
    @{ Cell = "E39"; Value = '  +0.48%  ' }
    @{ Cell = "D40"; Value = '6.272' }
    @{ Cell = "E40"; Value = '  +0.21%  ' }
    @{ Cell = "D41"; Value = '73.43' }
    @{ Cell = "E41"; Value = '  +6.00%  ' }
    @{ Cell = "D42"; Value = '1.995' }
    @{ Cell = "E42"; Value = '  +4.81%  ' }
    @{ Cell = "D43"; Value = '0.4243' }
    @{ Cell = "E43"; Value = '  +2.06%  ' }
    @{ Cell = "E44"; Value = '  +0.06%  ' }
    @{ Cell = "D45"; Value = '0.8380' }
    @{ Cell = "E45"; Value = '  -0.55%  ' }
    @{ Cell = "D46"; Value = '102.17' }
    @{ Cell = "D47"; Value = '9.369' }
    @{ Cell = "E47"; Value = '  +1.37%  ' }
    @{ Cell = "D48"; Value = '7.037' }
    @{ Cell = "E48"; Value = '  -0.53%  ' }
    @{ Cell = "D49"; Value = '35.47' }
    @{ Cell = "E49"; Value = '  +1.05%  ' }
    @{ Cell = "D50"; Value = '910.83' }
    @{ Cell = "E50"; Value = '  -0.09%  ' }
    @{ Cell = "B51"; Value = 'Decentraland' }
    @{ Cell = "C51"; Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana' }
    @{ Cell = "D51"; Value = '0.3917' }
    @{ Cell = "E51"; Value = '  +2.30%  ' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $isNumericLooking = $u.Value -match "^[+-]?[0-9]*\.?[0-9]+$"
    if ($isNumericLooking) {
        # Force text storage so Excel does not coerce the literal into a
        # real number (which would normalize "248.20" -> 248.2, etc.), then
        # restore the default style so no stray NumberFormat is left behind.
        $cell.NumberFormat = "@"
        $cell.Value = $u.Value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $u.Value
    }
}

Write-Output "updated $($updates.Count) cells"
